# Sync attendance_reports: fix "Recorded By" (column G) ordering.
#
# Each G cell holds a comma-separated list of recorder names
# (e.g. "dnasr281@gmail.com, System"). This edit rotates that list
# one position to the right, i.e. the last name in the list moves to
# the front: "a, b, c" -> "c, a, b". Rows whose G value has only a
# single name are unaffected by construction (rotating a 1-item list
# is a no-op).
#
# Only the specific rows touched by the upstream sync are updated here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$col = 7  # column G = "Recorded By"

$targetRows = @(
    2,3,4,5,6,7,8,10,12,13,14,15,18,19,20,21,22,24,26,
    28,29,30,31,32,33,34,36,38,39,40,41,
    44,45,46,47,48,50,52,
    54,55,56,57,58,59,60,62,64,65,66,67,
    70,71,72,73,74,76,78,
    80,81,82,83,84,85,86,87,
    90,92,
    99,101,
    106,107,108,109,110,111,112,113,
    116,118,
    125,127,
    132,133,134,135,136,137,138,139,
    142,144,
    151,153
)

function Rotate-Right-CsvList($value) {
    $parts = $value -split ","
    $parts = @($parts | ForEach-Object { $_.Trim() })

    if ($parts.Count -le 1) {
        return $value
    }

    $lastItem = $parts[$parts.Count - 1]
    $remaining = @($parts[0..($parts.Count - 2)])
    $rotated = @($lastItem) + $remaining

    return ($rotated -join ", ")
}

$updated = 0
foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, $col)
    $current = [string]$cell.Value2
    $newValue = Rotate-Right-CsvList $current
    if ($newValue -ne $current) {
        $cell.Value = $newValue
        $updated++
    }
}

Write-Host "Rotated Recorded By (column G) for $updated of $($targetRows.Count) target rows"
